$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.545.02'
$ws.Range('E2').Value = '  +1.43%  '

$ws.Range('D3').Value = '2.490.81'
$ws.Range('E3').Value = '  +1.86%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.995'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  -0.37%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '313.84'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +1.27%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '93.67'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.49%  '

$ws.Range('E7').Value = '  +0.04%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.996'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -0.41%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.498'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -0.74%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '32.95'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -1.09%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0785'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +0.76%  '

$ws.Range('E12').Value = '  +2.15%  '

$ws.Range('D13').Value = '2.875.79'
$ws.Range('E13').Value = '  +2.32%  '

$ws.Range('E14').Value = '  -0.61%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '15.57'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +8.31%  '

$ws.Range('D16').Value = '2.444.97'
$ws.Range('E16').Value = '  +1.18%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.760'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -2.85%  '

$ws.Range('D18').Value = '41.752.27'
$ws.Range('E18').Value = '  +2.02%  '

$ws.Range('E19').Value = '  +0.35%  '

$ws.Range('D20').Value = '0.0₃0923'
$ws.Range('E20').Value = '  +1.41%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '70.74'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +5.40%  '

$ws.Range('E22').Value = '  -1.86%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '236.40'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -0.02%  '

$ws.Range('E24').Value = '  -1.37%  '

$ws.Range('E25').Value = '  -0.89%  '

$ws.Range('E26').Value = '  -0.09%  '

$ws.Range('E27').Value = '  +1.35%  '

$ws.Range('E28').Value = '  +0.93%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '9.68'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +0.30%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '36.17'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +0.51%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '154.46'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +1.09%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '5.41'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -2.80%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '18.35'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +6.65%  '

$ws.Range('E34').Value = '  -0.99%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.0756'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +0.93%  '

$ws.Range('E36').Value = '  -1.50%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.95'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -1.84%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.83'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -2.65%  '

$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.103'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +0.00%  '

$ws.Range('B40').Value = 'Stellar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.113'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +0.12%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '4.14'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -0.15%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.998'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -0.38%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '19.90'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -5.93%  '

$ws.Range('D44').Value = '1.950.49'

$ws.Range('E45').Value = '  +0.69%  '

$ws.Range('E46').Value = '  -1.58%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '8.85'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +2.25%  '

$ws.Range('D48').Value = '2.734.57'
$ws.Range('E48').Value = '  +2.06%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '96.52'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -0.28%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '67.49'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -1.87%  '

$ws.Range('E51').Value = '  -2.09%  '
